# Update "想去人数" (interest count) figures on the 展览 (Exhibition),
# 本地生活 (Local Life) and 全部类型 (All types) sheets to match the
# newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) -------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 9071
$ws1.Range("F8").Value  = 6297
$ws1.Range("F12").Value = 9335
$ws1.Range("F13").Value = 10842
$ws1.Range("F14").Value = 1218
$ws1.Range("F15").Value = 1111
$ws1.Range("F16").Value = 4867
$ws1.Range("F17").Value = 780
$ws1.Range("F22").Value = 1319
$ws1.Range("F25").Value = 863
$ws1.Range("F29").Value = 411
$ws1.Range("F30").Value = 597
$ws1.Range("F31").Value = 2606
$ws1.Range("F33").Value = 178
$ws1.Range("F34").Value = 1698
$ws1.Range("F38").Value = 16
$ws1.Range("F39").Value = 902
$ws1.Range("F40").Value = 574
$ws1.Range("F41").Value = 3268
$ws1.Range("F45").Value = 568
$ws1.Range("F48").Value = 231
$ws1.Range("F49").Value = 4191

# --- Sheet: 本地生活 (Local Life) ----------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5788

# --- Sheet: 全部类型 (All types) ------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 9071
$ws4.Range("F8").Value  = 6297
$ws4.Range("F10").Value = 9335
$ws4.Range("F11").Value = 9335
$ws4.Range("F12").Value = 10842
$ws4.Range("F14").Value = 1218
$ws4.Range("F15").Value = 1111
$ws4.Range("F16").Value = 4867
$ws4.Range("F17").Value = 780
$ws4.Range("F22").Value = 1319
$ws4.Range("F24").Value = 863
$ws4.Range("F29").Value = 411
$ws4.Range("F30").Value = 2606
$ws4.Range("F31").Value = 178
$ws4.Range("F32").Value = 1698
$ws4.Range("F39").Value = 902
$ws4.Range("F40").Value = 574
$ws4.Range("F45").Value = 568
$ws4.Range("F47").Value = 231
$ws4.Range("F48").Value = 4191
